$wb = $excel.ActiveWorkbook

# Sheet "展览" (1st sheet): F2 104 -> 105, F3 306 -> 307
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 105
$ws1.Range("F3").Value = 307

# Sheet "全部类型" (4th sheet): F2 104 -> 105, F3 306 -> 307
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 105
$ws4.Range("F3").Value = 307
